# Updates the cryptocurrency price/volume table to reflect the latest
# GitHub Actions scrape (commit message: "Updated cryptos list ... with
# GitHub Actions"). Cell D (Price) and E (Volume(1h)) hold plain text in
# the source data, even though some of the new values look numeric
# (e.g. "211.81"). To avoid Excel's automatic number coercion (and the
# style/number-format changes that coercion would introduce), price
# values that parse as numbers are written through Set-TextValue, which
# forces a Text number format for the assignment and then restores the
# cell to the default ("Normal") style so the resulting XML has no
# spurious style attribute - matching the original file's plain,
# unstyled data cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue {
    param($cellRef, $text)
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

$ws.Range("D2").Value = '26.268.11'
$ws.Range("E2").Value = '  +0.18%  '
$ws.Range("D3").Value = '1.597.39'
$ws.Range("E3").Value = '  +0.69%  '
$ws.Range("E4").Value = '  +0.02%  '
Set-TextValue 'D5' '211.81'
$ws.Range("E5").Value = '  +0.01%  '
$ws.Range("E6").Value = '  +0.00%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").Value = '  +0.09%  '
$ws.Range("E9").Value = '  +0.38%  '
Set-TextValue 'D10' '18.99'
$ws.Range("E10").Value = '  -1.07%  '
$ws.Range("E11").Value = '  +0.77%  '
$ws.Range("E12").Value = '  +0.61%  '
$ws.Range("D13").Value = '1.585.71'
$ws.Range("E13").Value = '  +0.21%  '
$ws.Range("E14").Value = '  -0.27%  '
$ws.Range("E15").Value = '  -2.42%  '
Set-TextValue 'D16' '63.70'
$ws.Range("E16").Value = '  -0.34%  '
$ws.Range("D17").Value = '26.260.88'
$ws.Range("E17").Value = '  +0.18%  '
Set-TextValue 'D18' '230.32'
$ws.Range("E18").Value = '  +7.70%  '
Set-TextValue 'D19' '7.67'
$ws.Range("E19").Value = '  +4.80%  '
$ws.Range("D20").Value = '0.0₃0722'
$ws.Range("E20").Value = '  -0.28%  '
$ws.Range("E21").Value = '  +0.00%  '
$ws.Range("E22").Value = '  -0.17%  '
$ws.Range("B23").Value = 'Avalanche'
$ws.Range("C23").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue 'D23' '8.94'
$ws.Range("E23").Value = '  -0.35%  '
$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D24' '2.15'
$ws.Range("E24").Value = '  +1.11%  '
Set-TextValue 'D25' '145.93'
$ws.Range("E25").Value = '  +0.92%  '
$ws.Range("E26").Value = '  +0.04%  '
$ws.Range("E27").Value = '  +0.67%  '
$ws.Range("E28").Value = '  +0.39%  '
Set-TextValue 'D29' '15.33'
$ws.Range("E29").Value = '  +1.59%  '
$ws.Range("E30").Value = '  -0.44%  '
$ws.Range("E31").Value = '  +0.37%  '
Set-TextValue 'D32' '3.20'
$ws.Range("E32").Value = '  +0.82%  '
$ws.Range("D33").Value = '1.469.19'
$ws.Range("E33").Value = '  +4.11%  '
$ws.Range("E34").Value = '  +0.32%  '
$ws.Range("E35").Value = '  -0.38%  '
$ws.Range("E36").Value = '  +0.86%  '
Set-TextValue 'D37' '0.568'
$ws.Range("E37").Value = '  -3.11%  '
$ws.Range("E38").Value = '  -0.91%  '
Set-TextValue 'D39' '0.820'
$ws.Range("E39").Value = '  -0.04%  '
$ws.Range("E40").Value = '  -2.06%  '
$ws.Range("E41").Value = '  +0.06%  '
Set-TextValue 'D42' '2.18'
$ws.Range("E42").Value = '  +2.31%  '
Set-TextValue 'D43' '0.933'
$ws.Range("E43").Value = '  -2.05%  '
$ws.Range("D44").Value = '1.734.04'
$ws.Range("E44").Value = '  +0.79%  '
Set-TextValue 'D45' '0.756'
$ws.Range("E45").Value = '  -1.34%  '
Set-TextValue 'D46' '60.60'
$ws.Range("E46").Value = '  -0.60%  '
$ws.Range("E47").Value = '  +2.72%  '
$ws.Range("E48").Value = '  -0.38%  '
$ws.Range("E49").Value = '  +0.07%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D50' '7.43'
$ws.Range("E50").Value = '  +1.35%  '
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 'D51' '0.0948'
$ws.Range("E51").Value = '  -2.02%  '
